{"js": "// The generated Java stack trace in the document's second paragraph was\n// regenerated against a newer build (line numbers / frames shifted, e.g.\n// AstEvaluator.java:186 -> 189, M2DocEvaluator.java:1343 -> 1462, plus the\n// JUnit internals were inlined differently by a newer JUnit/JDK).\n// We locate the run holding that stack trace via its distinctive first and\n// last lines, then swap its text for the updated trace, keeping the run's\n// bold/red formatting and the trailing manual line break untouched.\nconst body = context.document.body;\n\nconst startMatches = body.search(\"divOp(java.lang.Integer\", { matchCase: true });\nstartMatches.load(\"items\");\nconst endMatches = body.search(\"RemoteTestRunner.java:206)\", { matchCase: true });\nendMatches.load(\"items\");\nawait context.sync();\n\nif (startMatches.items.length === 0 || endMatches.items.length === 0) {\n  throw new Error(\"Could not locate the stack trace text to update.\");\n}\n\nconst traceRange = startMatches.items[0].expandTo(endMatches.items[0]);\n\nconst newTrace = \"divOp(java.lang.Integer,java.lang.Integer) with arguments [1, 0] failed:\\n\\t/ by zero\\njava.lang.ArithmeticException: / by zero\\n\\tat org.eclipse.acceleo.query.services.NumberServices.divOp(NumberServices.java:99)\\n\\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\\n\\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\\n\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\\n\\tat java.lang.reflect.Method.invoke(Method.java:498)\\n\\tat org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:162)\\n\\tat org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:135)\\n\\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callService(EvaluationServices.java:129)\\n\\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:172)\\n\\tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:189)\\n\\tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:119)\\n\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\\n\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\\n\\tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:112)\\n\\tat org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:52)\\n\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseLet(M2DocEvaluator.java:1462)\\n\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseLet(M2DocEvaluator.java:1)\\n\\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:314)\\n\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\\n\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\\n\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1450)\\n\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1675)\\n\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)\\n\\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)\\n\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\\n\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\\n\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1450)\\n\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:314)\\n\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)\\n\\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)\\n\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\\n\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\\n\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1450)\\n\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:299)\\n\\tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:853)\\n\\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:536)\\n\\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:422)\\n\\tat sun.reflect.GeneratedMethodAccessor6.invoke(Unknown Source)\\n\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\\n\\tat java.lang.reflect.Method.invoke(Method.java:498)\\n\\tat org.junit.runners.model.FrameworkMethod$1.runReflectiveCall(FrameworkMethod.java:50)\\n\\tat org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12)\\n\\tat org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:47)\\n\\tat org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)\\n\\tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)\\n\\tat org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:325)\\n\\tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:78)\\n\\tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:57)\\n\\tat org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\\n\\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\\n\\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\\n\\tat org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\\n\\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\\n\\tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)\\n\\tat org.junit.runners.Suite.runChild(Suite.java:128)\\n\\tat org.junit.runners.Suite.runChild(Suite.java:27)\\n\\tat org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\\n\\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\\n\\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\\n\\tat org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\\n\\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\\n\\tat org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)\\n\\tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)\\n\\tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)\\n\\tat org.junit.runners.Suite.runChild(Suite.java:128)\\n\\tat org.junit.runners.Suite.runChild(Suite.java:27)\\n\\tat org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\\n\\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\\n\\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\\n\\tat org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\\n\\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\\n\\tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)\\n\\tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)\\n\\tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)\\n\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)\\n\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)\\n\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)\\n\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)\";\n\ntraceRange.insertText(newTrace, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The generated Java stack trace in the document's second paragraph was\n# regenerated against a newer build (line numbers / frames shifted, e.g.\n# AstEvaluator.java:186 -> 189, M2DocEvaluator.java:1343 -> 1462, plus the\n# JUnit internals were inlined differently by a newer JUnit/JDK).\n# We locate the run holding that stack trace by searching for its\n# distinctive first and last lines, then swap its text for the updated\n# trace, keeping the run's bold/red formatting and the trailing manual\n# line break untouched.\n$d = $word.ActiveDocument\n\n$startRange = $d.Content\n$startFind = $startRange.Find\n$startFound = $startFind.Execute(\"divOp(java.lang.Integer\")\nif (-not $startFound) {\n    throw \"Could not find the start of the stack trace.\"\n}\n$startPos = $startRange.Start\n\n$endRange = $d.Content\n$endFind = $endRange.Find\n$endFound = $endFind.Execute(\"RemoteTestRunner.java:206)\")\nif (-not $endFound) {\n    throw \"Could not find the end of the stack trace.\"\n}\n$endPos = $endRange.End\n\n$traceRange = $d.Range($startPos, $endPos)\n\n$newTrace = 'divOp(java.lang.Integer,java.lang.Integer) with arguments [1, 0] failed:\u0001NL\u0001\u0001TAB\u0001/ by zero\u0001NL\u0001java.lang.ArithmeticException: / by zero\u0001NL\u0001\u0001TAB\u0001at org.eclipse.acceleo.query.services.NumberServices.divOp(NumberServices.java:99)\u0001NL\u0001\u0001TAB\u0001at sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\u0001NL\u0001\u0001TAB\u0001at sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\u0001NL\u0001\u0001TAB\u0001at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\u0001NL\u0001\u0001TAB\u0001at java.lang.reflect.Method.invoke(Method.java:498)\u0001NL\u0001\u0001TAB\u0001at org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:162)\u0001NL\u0001\u0001TAB\u0001at org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:135)\u0001NL\u0001\u0001TAB\u0001at org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callService(EvaluationServices.java:129)\u0001NL\u0001\u0001TAB\u0001at org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:172)\u0001NL\u0001\u0001TAB\u0001at org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:189)\u0001NL\u0001\u0001TAB\u0001at org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:119)\u0001NL\u0001\u0001TAB\u0001at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\u0001NL\u0001\u0001TAB\u0001at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\u0001NL\u0001\u0001TAB\u0001at org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:112)\u0001NL\u0001\u0001TAB\u0001at org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:52)\u0001NL\u0001\u0001TAB\u0001at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseLet(M2DocEvaluator.java:1462)\u0001NL\u0001\u0001TAB\u0001at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseLet(M2DocEvaluator.java:1)\u0001NL\u0001\u0001TAB\u0001at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:314)\u0001NL\u0001\u0001TAB\u0001at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\u0001NL\u0001\u0001TAB\u0001at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\u0001NL\u0001\u0001TAB\u0001at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1450)\u0001NL\u0001\u0001TAB\u0001at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1675)\u0001NL\u0001\u0001TAB\u0001at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)\u0001NL\u0001\u0001TAB\u0001at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)\u0001NL\u0001\u0001TAB\u0001at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\u0001NL\u0001\u0001TAB\u0001at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\u0001NL\u0001\u0001TAB\u0001at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1450)\u0001NL\u0001\u0001TAB\u0001at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:314)\u0001NL\u0001\u0001TAB\u0001at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)\u0001NL\u0001\u0001TAB\u0001at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)\u0001NL\u0001\u0001TAB\u0001at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\u0001NL\u0001\u0001TAB\u0001at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\u0001NL\u0001\u0001TAB\u0001at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1450)\u0001NL\u0001\u0001TAB\u0001at org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:299)\u0001NL\u0001\u0001TAB\u0001at org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:853)\u0001NL\u0001\u0001TAB\u0001at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:536)\u0001NL\u0001\u0001TAB\u0001at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:422)\u0001NL\u0001\u0001TAB\u0001at sun.reflect.GeneratedMethodAccessor6.invoke(Unknown Source)\u0001NL\u0001\u0001TAB\u0001at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\u0001NL\u0001\u0001TAB\u0001at java.lang.reflect.Method.invoke(Method.java:498)\u0001NL\u0001\u0001TAB\u0001at org.junit.runners.model.FrameworkMethod$1.runReflectiveCall(FrameworkMethod.java:50)\u0001NL\u0001\u0001TAB\u0001at org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12)\u0001NL\u0001\u0001TAB\u0001at org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:47)\u0001NL\u0001\u0001TAB\u0001at org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)\u0001NL\u0001\u0001TAB\u0001at org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)\u0001NL\u0001\u0001TAB\u0001at org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:325)\u0001NL\u0001\u0001TAB\u0001at org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:78)\u0001NL\u0001\u0001TAB\u0001at org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:57)\u0001NL\u0001\u0001TAB\u0001at org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\u0001NL\u0001\u0001TAB\u0001at org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\u0001NL\u0001\u0001TAB\u0001at org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\u0001NL\u0001\u0001TAB\u0001at org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\u0001NL\u0001\u0001TAB\u0001at org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\u0001NL\u0001\u0001TAB\u0001at org.junit.runners.ParentRunner.run(ParentRunner.java:363)\u0001NL\u0001\u0001TAB\u0001at org.junit.runners.Suite.runChild(Suite.java:128)\u0001NL\u0001\u0001TAB\u0001at org.junit.runners.Suite.runChild(Suite.java:27)\u0001NL\u0001\u0001TAB\u0001at org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\u0001NL\u0001\u0001TAB\u0001at org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\u0001NL\u0001\u0001TAB\u0001at org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\u0001NL\u0001\u0001TAB\u0001at org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\u0001NL\u0001\u0001TAB\u0001at org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\u0001NL\u0001\u0001TAB\u0001at org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)\u0001NL\u0001\u0001TAB\u0001at org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)\u0001NL\u0001\u0001TAB\u0001at org.junit.runners.ParentRunner.run(ParentRunner.java:363)\u0001NL\u0001\u0001TAB\u0001at org.junit.runners.Suite.runChild(Suite.java:128)\u0001NL\u0001\u0001TAB\u0001at org.junit.runners.Suite.runChild(Suite.java:27)\u0001NL\u0001\u0001TAB\u0001at org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\u0001NL\u0001\u0001TAB\u0001at org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\u0001NL\u0001\u0001TAB\u0001at org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\u0001NL\u0001\u0001TAB\u0001at org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\u0001NL\u0001\u0001TAB\u0001at org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\u0001NL\u0001\u0001TAB\u0001at org.junit.runners.ParentRunner.run(ParentRunner.java:363)\u0001NL\u0001\u0001TAB\u0001at org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)\u0001NL\u0001\u0001TAB\u0001at org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)\u0001NL\u0001\u0001TAB\u0001at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)\u0001NL\u0001\u0001TAB\u0001at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)\u0001NL\u0001\u0001TAB\u0001at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)\u0001NL\u0001\u0001TAB\u0001at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)'\n$newTrace = $newTrace.Replace([char]1 + 'TAB' + [char]1, [string][char]9)\n$newTrace = $newTrace.Replace([char]1 + 'NL' + [char]1, [string][char]10)\n\n$traceRange.Text = $newTrace\n"}
